$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

$updates = @{
    "E2"  = 62
    "E9"  = 32
    "F9"  = 16
    "H9"  = 25
    "E12" = 10
    "F12" = 4
    "H12" = 4
    "E15" = 179
    "F15" = 100
    "H15" = 141
    "E25" = 26
    "F25" = 14
    "H25" = 22
    "E26" = 34
    "F26" = 18
    "H26" = 28
    "E27" = 18
    "E36" = 123
    "F36" = 59
    "H36" = 91
    "E47" = 64
    "F47" = 42
    "H47" = 52
    "E49" = 79
    "E60" = 23
    "F60" = 12
    "H60" = 17
    "E66" = 39
    "F66" = 28
    "H66" = 36
    "F67" = 26
    "G67" = 8
    "F70" = 26
    "G70" = 13
    "E76" = 59
    "E79" = 47
    "F79" = 23
    "H79" = 34
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
